$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.101.21"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "'3.765.69"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'623.76"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'180.04"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'3.764.21"
$ws.Range("E7").Value = "  +3.13%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("D11").Value = "'6.33"
$ws.Range("E11").Value = "  -5.03%  "
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "'41.25"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "'0.0000261"
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("D15").Value = "'4.384.20"
$ws.Range("E15").Value = "  +2.78%  "
$ws.Range("D16").Value = "'3.757.21"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "'70.168.63"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'7.67"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "'16.79"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "'507.93"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'87.26"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "'13.23"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "'11.15"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "'0.0000138"
$ws.Range("E28").Value = "  +25.72%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'2.96"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "'7.92"
$ws.Range("E32").Value = "  -2.84%  "
$ws.Range("D33").Value = "'31.55"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "'0.996"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'1.07"
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("D37").Value = "'6.25"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").Value = "'0.337"
$ws.Range("E38").Value = "  -3.32%  "
$ws.Range("D39").Value = "'0.133"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").Value = "'50.36"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "'45.24"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "'424.91"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").Value = "'8.77"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "'3.008.73"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").Value = "'0.0366"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").Value = "'27.45"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("D49").Value = "'138.82"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'2.55"
$ws.Range("E51").Value = "  +3.24%  "
